# chore: update Sheets via scheduled runner
# Refreshes currentAveragePrice/NQ/HQ, LevePriceNQ/HQ and recomputed
# LeveProfitNQ/HQ figures across the per-job Leve profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H95").Value = 49812
$ws.Range("J95").Value = 49812
$ws.Range("L95").Value = 49812
$ws.Range("N95").Value = -55304
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H137").Value = 5000
$ws.Range("I137").Value = 5000
$ws.Range("K137").Value = 15000
$ws.Range("M137").Value = -12450
$ws.Range("H138").Value = 1310
$ws.Range("I138").Value = 583.3333
$ws.Range("J138").Value = 2400
$ws.Range("K138").Value = 1749.9999
$ws.Range("L138").Value = 7200
$ws.Range("M138").Value = 3390.0001
$ws.Range("N138").Value = -17480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 56
$ws.Range("I5").Value = 56
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 56
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 56
$ws.Range("N5").ClearContents()
$ws.Range("H95").Value = 20653.75
$ws.Range("J95").Value = 20653.75
$ws.Range("L95").Value = 20653.75
$ws.Range("N95").Value = -26145.75
$ws.Range("H96").Value = 32468.8
$ws.Range("J96").Value = 32468.8
$ws.Range("L96").Value = 32468.8
$ws.Range("N96").Value = -37960.8
$ws.Range("H97").Value = 1749.75
$ws.Range("I97").Value = 1833
$ws.Range("K97").Value = 1833
$ws.Range("M97").Value = -1337
$ws.Range("H132").Value = 5369.364
$ws.Range("I132").Value = 3477.4285
$ws.Range("J132").Value = 8680.25
$ws.Range("K132").Value = 10432.2855
$ws.Range("L132").Value = 26040.75
$ws.Range("M132").Value = -7902.2855
$ws.Range("N132").Value = -31100.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 56
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 56
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 59
$ws.Range("N4").ClearContents()
$ws.Range("H36").Value = 8100
$ws.Range("I36").Value = 6375
$ws.Range("K36").Value = 6375
$ws.Range("M36").Value = -5841
$ws.Range("H107").Value = 2312.5
$ws.Range("I107").Value = 2083.3333
$ws.Range("K107").Value = 2083.3333
$ws.Range("M107").Value = -163.3332999999998
$ws.Range("H134").Value = 7086.231
$ws.Range("I134").Value = 3668.5
$ws.Range("K134").Value = 11005.5
$ws.Range("M134").Value = -8470.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3430
$ws.Range("I16").Value = 3430
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3430
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3143
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 4529.1665
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 4529.1665
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H58").Value = 1499.5
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 1999
$ws.Range("K58").Value = 1000
$ws.Range("L58").Value = 1999
$ws.Range("M58").Value = -797
$ws.Range("N58").Value = -2405
$ws.Range("H86").Value = 3656.8572
$ws.Range("I86").Value = 2959.6
$ws.Range("K86").Value = 2959.6
$ws.Range("M86").Value = -1836.6
$ws.Range("H89").Value = 3656.8572
$ws.Range("I89").Value = 2959.6
$ws.Range("K89").Value = 14798
$ws.Range("M89").Value = -9182
$ws.Range("H113").Value = 3430
$ws.Range("I113").Value = 3430
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3430
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1260
$ws.Range("N113").ClearContents()
$ws.Range("H125").Value = 34808.668
$ws.Range("J125").Value = 34808.668
$ws.Range("L125").Value = 34808.668
$ws.Range("N125").Value = -39728.668
$ws.Range("H132").Value = 2999
$ws.Range("I132").Value = 2499
$ws.Range("J132").Value = 3199
$ws.Range("K132").Value = 7497
$ws.Range("L132").Value = 9597
$ws.Range("M132").Value = -4967
$ws.Range("N132").Value = -14657
$ws.Range("H134").Value = 1837.3334
$ws.Range("I134").Value = 1837.3334
$ws.Range("K134").Value = 5512.0002
$ws.Range("M134").Value = -2977.0002
$ws.Range("H136").Value = 1499.5
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -11097

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 293.6
$ws.Range("I23").Value = 235.5
$ws.Range("J23").Value = 332.33334
$ws.Range("K23").Value = 706.5
$ws.Range("L23").Value = 997.0000200000001
$ws.Range("M23").Value = -471.5
$ws.Range("N23").Value = -1467.00002
$ws.Range("H141").Value = 3209.5715
$ws.Range("I141").Value = 3216.75
$ws.Range("K141").Value = 9650.25
$ws.Range("M141").Value = -4470.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 90000
$ws.Range("I93").Value = 90000
$ws.Range("K93").Value = 90000
$ws.Range("M93").Value = -88128
$ws.Range("H107").Value = 875
$ws.Range("I107").Value = 875
$ws.Range("K107").Value = 875
$ws.Range("M107").Value = 1045
$ws.Range("H132").Value = 7550.5
$ws.Range("I132").Value = 3990
$ws.Range("K132").Value = 11970
$ws.Range("M132").Value = -9440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1716.5
$ws.Range("I93").Value = 1716.5
$ws.Range("K93").Value = 1716.5
$ws.Range("M93").Value = -468.5
$ws.Range("H125").Value = 57283.5
$ws.Range("J125").Value = 57283.5
$ws.Range("L125").Value = 57283.5
$ws.Range("N125").Value = -67123.5
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
$ws.Range("H136").Value = 82958
$ws.Range("I136").Value = 3559.8
$ws.Range("K136").Value = 10679.4
$ws.Range("M136").Value = -8129.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 35000.5
$ws.Range("J70").Value = 40000.332
$ws.Range("L70").Value = 40000.332
$ws.Range("N70").Value = -40630.332
$ws.Range("H73").Value = 35000.5
$ws.Range("J73").Value = 40000.332
$ws.Range("L73").Value = 40000.332
$ws.Range("N73").Value = -42184.332
$ws.Range("H100").Value = 687.25
$ws.Range("I100").Value = 687.25
$ws.Range("K100").Value = 1374.5
$ws.Range("M100").Value = -833.5
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
